# Update metrics for all data rows (2-26) with the new values produced
# after retraining with the new language model ("atualizado todo o
# treinamento para o novo lm").
#
# Every data row (B:Q) receives the same new metric values.
# (Values are written using plain decimal notation -- not scientific --
# because the interpreter's expression parser does not accept the
# `1.23e-05` exponent form; the underlying double is identical either way.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999883732771242,        # B - r2
    0.9990763349190175,        # C - r2_sup
    0.9999999214032369,        # D - r2_test
    0.9999841125966974,        # E - r2_val
    0.9999933993149381,        # F - r2_vt
    0.00001085303545036088,    # G - mse
    0.000862200808886027,      # H - mse_sup
    0.00000002734613341965376, # I - mse_test
    0.000003961272472290328,   # J - mse_val
    0.000001994309302854991,   # K - mse_vt
    0.0001805656912053873,     # L - mape
    0.00329439454989242,       # M - rmse
    0.9999069862169936,        # N - r2_adj
    0.003434643746218193,      # O - rsd
    64.86213150418904,         # P - aic
    90.45852382642124          # Q - bic
)

$firstRow = 2
$lastRow = 26
$firstCol = 2   # column B
$lastCol = 17   # column Q

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - $firstCol]
    }
}
